$d = $word.ActiveDocument

$replacements = @(
    @{old="933÷9="; new="958÷5="},
    @{old="443÷9="; new="128÷9="},
    @{old="700÷6="; new="458÷8="},
    @{old="579÷9="; new="943÷9="},
    @{old="866÷5="; new="564÷4="},
    @{old="825÷6="; new="717÷6="},
    @{old="309÷6="; new="436÷3="},
    @{old="410÷8="; new="525÷8="},
    @{old="154÷5="; new="326÷5="},
    @{old="472÷8="; new="239÷5="},
    @{old="763÷4="; new="787÷9="},
    @{old="738÷2="; new="199÷7="},
    @{old="355÷5="; new="401÷2="},
    @{old="321÷8="; new="538÷9="},
    @{old="970÷8="; new="548÷9="},
    @{old="478÷3="; new="355÷5="},
    @{old="956÷7="; new="606÷9="},
    @{old="978÷8="; new="706÷7="},
    @{old="788÷3="; new="752÷6="},
    @{old="830÷4="; new="850÷5="},
    @{old="782÷8="; new="607÷8="},
    @{old="643÷9="; new="458÷2="},
    @{old="136÷9="; new="961÷2="},
    @{old="451÷7="; new="662÷9="},
    @{old="878÷9="; new="244÷4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
